$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 151
$ws.Range("I2").Value = 347
$ws.Range("J2").Value = 1236
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 351
$ws.Range("M2").Value = 26
$ws.Range("N2").Value = 232
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 13
$ws.Range("S2").Value = 152
$ws.Range("T2").Value = 238
$ws.Range("U2").Value = 12
$ws.Range("V2").Value = 2041
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 2074
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 30
$ws.Range("AA2").Value = 8
